# Add Setups and Results
# Move the Mean/Standard Deviation values out of columns B/C and into the
# Max/Min columns (D/E) with updated values for rows 2-4 on the "DOE" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DOE")

# Clear the old B/C values (Mean / Standard Deviation)
$ws.Range("B2:C4").ClearContents()

# Row 2 (POX/C): Max = 1000, Min = 10
$ws.Range("D2").Value = 1000
$ws.Range("E2").Value = 10

# Row 3 (C/A): Max = 1, Min = 0.001
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.001

# Row 4 (POX/M): Max = 0.002, Min = 0.00083
$ws.Range("D4").Value = 0.002
$ws.Range("E4").Value = 0.00083
